# Add speaker notes to the first two slides and append a new blank slide
# (duplicate of the existing "Blank" layout slide), matching the commit
# that expands the deck with protocol-description notes pages plus a
# placeholder slide 3 for further work.

$p = $ppt.ActivePresentation

# Slide 1: "Illustrations" title slide -> notes "Joystick Experimental Setup"
$slide1 = $p.Slides.Item(1)
$notes1 = $slide1.NotesPage.Shapes.Placeholders.Item(2)
$notes1.TextFrame.TextRange.Text = "Joystick Experimental Setup"

# Slide 2: blank slide -> notes "LabBench I/O Experimental Setup"
$slide2 = $p.Slides.Item(2)
$notes2 = $slide2.NotesPage.Shapes.Placeholders.Item(2)
$notes2.TextFrame.TextRange.Text = "LabBench I/O Experimental Setup"

# Append a new slide 3 by duplicating the blank slide 2 (keeps the same
# "Blank" layout and stays free of any shapes/notes of its own).
$slide2.Duplicate() | Out-Null

Write-Output ("Slides: " + $p.Slides.Count)
